# Daily attendance processing - 2025-10-06 13:29:23
# Rotate the "Recorded By" (column G) list of names/emails for the specific
# rows that were re-processed: move the first comma-separated entry to the
# end of the list (left-rotate by one position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,11,12,13,29,30,32,33,38,39,40,56,57,58,59,60,65,66,67,84,85,89,90,93,110,111,115,116,119,136,137,141,142,145)

foreach ($r in $rows) {
    $cell = $ws.Range("G$r")
    $text = $cell.Text
    $parts = $text -split ", "
    if ($parts.Length -gt 1) {
        $rotated = ($parts[1..($parts.Length - 1)] + $parts[0]) -join ", "
        $cell.Value = $rotated
    }
}
